$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 1 (the old numeric A1 cell); this shifts the former row 2
# (shared-string cell) up to become the new row 1 / A1.
$ws.Rows("1").Delete()

# Replace the shared string text with the reformatted (pretty-printed
# JSON-like) content.
$text = @'
questions = [
    {
        "title": "You are applying steepest descent hill-climbing for allocating electric pillars in a given region. You have the following values of objective function for the four next possible moves. The current state\u2019s objective function value is 8. Which move should you pick?",
        "ques_type": 2,
        "options": [
            "1",
            "2",
            "3",
            "4"
        ],
        "score": "2"
    },
    {
        "title": "You are programming a game to play tic-tac-toe using adversarial search. You want to check whether the current game configuration is endpoint. Which function from the formal definition of games should you use?",
        "ques_type": 2,
        "options": [
            "Actions()",
            "Result()",
            "Is-Terminal()",
            "Utility()"
        ],
        "score": "Is-Terminal()"
    },
    {
        "title": "The problem you are investigating contains eight possible models. Knowledge base (KB) is true in five of those models. You want to check whether KB entails sentence \u0251. In how many of the models where KB is true must \u0251 itself be true to be entailed from KB?",
        "ques_type": 2,
        "options": [
            "1",
            "3",
            "5",
            "8"
        ],
        "score": "5"
    },
    {
        "title": "You are designing a planning domain definition language problem and have the following fluents: F1, F2, and F3 Which of the following definitions describes the initial state of the problem?",
        "ques_type": 2,
        "options": [
            "F1 \u2228 F2 \u2228 F3",
            "F1 \u2227 F2 \u2228 F3",
            "F1 \u2228 F2 \u2227 F3",
            "F1 \u2227 F2 \u2227 F3"
        ],
        "score": "F1 \u2227 F2 \u2227 F3"
    }
]
'@
$ws.Range("A1").Value = $text
